$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A68").Value = 2
$ws.Range("B68").Value = "Devolve with Xerogen Crystals"
$ws.Range("C68").Value = "King Hulk, Sakaarson"
$ws.Range("D68").Value = "Warbound|Maximum Carnage"
$ws.Range("E68").Value = "Shi'ar Patrol Craft|The Brood"
$ws.Range("F68").Value = "Nerkkod, Breaker of Oceans (FI)|X-23 (XM)|Hulkling (CW)|Gambit (B)|Angel Noir (N)"
$ws.Range("G68").Value = 1
$ws.Range("J68").Value = "The Brood"
$ws.Range("H68").Value = "40|23"
$ws.Range("I68").Value = "not really"
$ws.Range("K68").Value = "Not many villain cards left, but only two escapees. Lots of wound clearers and top deck manipulation + feast thinned decks well."

$ws.Range("K69").Select()
